$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1120.8667
$ws.Range("J19").Value = 1054.7693
$ws.Range("L19").Value = 1054.7693
$ws.Range("N19").Value = -1404.7693
$ws.Range("H88").Value = 74915.73
$ws.Range("I88").Value = 4859.6665
$ws.Range("J88").Value = 101186.75
$ws.Range("K88").Value = 4859.6665
$ws.Range("L88").Value = 101186.75
$ws.Range("M88").Value = -4453.6665
$ws.Range("N88").Value = -101998.75
$ws.Range("H91").Value = 74915.73
$ws.Range("I91").Value = 4859.6665
$ws.Range("J91").Value = 101186.75
$ws.Range("K91").Value = 4859.6665
$ws.Range("L91").Value = 101186.75
$ws.Range("M91").Value = -3455.6665
$ws.Range("N91").Value = -103994.75
$ws.Range("H113").Value = 18531518
$ws.Range("J113").Value = 18500
$ws.Range("L113").Value = 18500
$ws.Range("N113").Value = -25008
$ws.Range("H132").Value = 1925.95
$ws.Range("I132").Value = 1974.6842
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 5924.0526
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -3394.0526
$ws.Range("N132").Value = -8060
$ws.Range("H137").Value = 3161.077
$ws.Range("I137").Value = 5057.6
$ws.Range("K137").Value = 15172.8
$ws.Range("M137").Value = -12622.8
$ws.Range("H138").Value = 5241.533
$ws.Range("I138").Value = 1525.1818
$ws.Range("J138").Value = 7393.1055
$ws.Range("K138").Value = 4575.5454
$ws.Range("L138").Value = 22179.3165
$ws.Range("M138").Value = 564.4546
$ws.Range("N138").Value = -32459.3165
$ws.Range("H141").Value = 1552.9546
$ws.Range("I141").Value = 1552.9546
$ws.Range("K141").Value = 4658.8638
$ws.Range("M141").Value = 521.1361999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3475542
$ws.Range("I32").Value = 3790818.2
$ws.Range("K32").Value = 3790818.2
$ws.Range("M32").Value = -3790531.2
$ws.Range("H76").Value = 44000
$ws.Range("J76").Value = 44000
$ws.Range("L76").Value = 44000
$ws.Range("N76").Value = -44676
$ws.Range("H79").Value = 44000
$ws.Range("J79").Value = 44000
$ws.Range("L79").Value = 44000
$ws.Range("N79").Value = -46340
$ws.Range("H88").Value = 1914.625
$ws.Range("I88").Value = 1289.091
$ws.Range("J88").Value = 2443.923
$ws.Range("K88").Value = 1289.091
$ws.Range("L88").Value = 2443.923
$ws.Range("M88").Value = -883.0909999999999
$ws.Range("N88").Value = -3255.923
$ws.Range("H91").Value = 1914.625
$ws.Range("I91").Value = 1289.091
$ws.Range("J91").Value = 2443.923
$ws.Range("K91").Value = 1289.091
$ws.Range("L91").Value = 2443.923
$ws.Range("M91").Value = 114.9090000000001
$ws.Range("N91").Value = -5251.923
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5326.8716
$ws.Range("I134").Value = 3454.451
$ws.Range("J134").Value = 10352.842
$ws.Range("K134").Value = 10363.353
$ws.Range("L134").Value = 31058.526
$ws.Range("M134").Value = -7828.352999999999
$ws.Range("N134").Value = -36128.526
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5866.533
$ws.Range("J16").Value = 7332.1113
$ws.Range("L16").Value = 7332.1113
$ws.Range("N16").Value = -7906.1113
$ws.Range("H31").Value = 8367.139999999999
$ws.Range("I31").Value = 3426.3
$ws.Range("K31").Value = 3426.3
$ws.Range("M31").Value = -3131.3
$ws.Range("H34").Value = 8367.139999999999
$ws.Range("I34").Value = 3426.3
$ws.Range("K34").Value = 3426.3
$ws.Range("M34").Value = -3224.3
$ws.Range("H47").Value = 36999
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 36999
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 36999
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -38131
$ws.Range("H58").Value = 8777430
$ws.Range("I58").Value = 17859466
$ws.Range("K58").Value = 17859466
$ws.Range("M58").Value = -17859263
$ws.Range("H62").Value = 9619984
$ws.Range("I62").Value = 20836934
$ws.Range("J62").Value = 5456.2856
$ws.Range("K62").Value = 20836934
$ws.Range("L62").Value = 5456.2856
$ws.Range("M62").Value = -20836310
$ws.Range("N62").Value = -6704.2856
$ws.Range("H65").Value = 9619984
$ws.Range("I65").Value = 20836934
$ws.Range("J65").Value = 5456.2856
$ws.Range("K65").Value = 104184670
$ws.Range("L65").Value = 27281.428
$ws.Range("M65").Value = -104181550
$ws.Range("N65").Value = -33521.428
$ws.Range("H99").Value = 12614
$ws.Range("I99").Value = 19637.334
$ws.Range("J99").Value = 8400
$ws.Range("K99").Value = 19637.334
$ws.Range("L99").Value = 8400
$ws.Range("M99").Value = -18139.334
$ws.Range("N99").Value = -11396
$ws.Range("H113").Value = 5866.533
$ws.Range("J113").Value = 7332.1113
$ws.Range("L113").Value = 7332.1113
$ws.Range("N113").Value = -11672.1113
$ws.Range("H126").Value = 12614
$ws.Range("I126").Value = 19637.334
$ws.Range("J126").Value = 8400
$ws.Range("K126").Value = 58912.00199999999
$ws.Range("L126").Value = 25200
$ws.Range("M126").Value = -56442.00199999999
$ws.Range("N126").Value = -30140
$ws.Range("H132").Value = 4830.9487
$ws.Range("I132").Value = 1763.7307
$ws.Range("J132").Value = 10965.385
$ws.Range("K132").Value = 5291.1921
$ws.Range("L132").Value = 32896.155
$ws.Range("M132").Value = -2761.1921
$ws.Range("N132").Value = -37956.155
$ws.Range("H134").Value = 6824.5
$ws.Range("I134").Value = 5536.5557
$ws.Range("J134").Value = 8932.046
$ws.Range("K134").Value = 16609.6671
$ws.Range("L134").Value = 26796.138
$ws.Range("M134").Value = -14074.6671
$ws.Range("N134").Value = -31866.138
$ws.Range("H136").Value = 8777430
$ws.Range("I136").Value = 17859466
$ws.Range("K136").Value = 53578398
$ws.Range("M136").Value = -53575848
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 647.86664
$ws.Range("I117").Value = 487.5
$ws.Range("J117").Value = 706.1818
$ws.Range("K117").Value = 1462.5
$ws.Range("L117").Value = 2118.5454
$ws.Range("M117").Value = 1979.5
$ws.Range("N117").Value = -9002.545399999999
$ws.Range("H122").Value = 2832126.5
$ws.Range("J122").Value = 7499.5
$ws.Range("L122").Value = 67495.5
$ws.Range("N122").Value = -72395.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 801446
$ws.Range("J107").Value = 1660
$ws.Range("L107").Value = 1660
$ws.Range("N107").Value = -5500
$ws.Range("H113").Value = 5795.147
$ws.Range("I113").Value = 2820.3076
$ws.Range("K113").Value = 2820.3076
$ws.Range("M113").Value = -650.3076000000001
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 3438.9487
$ws.Range("I132").Value = 3150.8484
$ws.Range("J132").Value = 5023.5
$ws.Range("K132").Value = 9452.5452
$ws.Range("L132").Value = 15070.5
$ws.Range("M132").Value = -6922.5452
$ws.Range("N132").Value = -20130.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5866.885
$ws.Range("I7").Value = 3854
$ws.Range("J7").Value = 7879.769
$ws.Range("K7").Value = 3854
$ws.Range("L7").Value = 7879.769
$ws.Range("M7").Value = -3742
$ws.Range("N7").Value = -8103.769
$ws.Range("H16").Value = 1892.9524
$ws.Range("I16").Value = 1787
$ws.Range("K16").Value = 1787
$ws.Range("M16").Value = -1617
$ws.Range("H121").Value = 56959
$ws.Range("J121").Value = 56959
$ws.Range("L121").Value = 56959
$ws.Range("N121").Value = -60453
$ws.Range("H122").Value = 5489.16
$ws.Range("I122").Value = 4643.3687
$ws.Range("K122").Value = 13930.1061
$ws.Range("M122").Value = -11480.1061
$ws.Range("H126").Value = 5866.885
$ws.Range("I126").Value = 3854
$ws.Range("J126").Value = 7879.769
$ws.Range("K126").Value = 11562
$ws.Range("L126").Value = 23639.307
$ws.Range("M126").Value = -9092
$ws.Range("N126").Value = -28579.307
$ws.Range("H132").Value = 15635310
$ws.Range("I132").Value = 55565264
$ws.Range("K132").Value = 166695792
$ws.Range("M132").Value = -166693262
$ws.Range("H136").Value = 12622.091
$ws.Range("I136").Value = 7896.6665
$ws.Range("J136").Value = 13368.211
$ws.Range("K136").Value = 23689.9995
$ws.Range("L136").Value = 40104.633
$ws.Range("M136").Value = -21139.9995
$ws.Range("N136").Value = -45204.633
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 755.1875
$ws.Range("I100").Value = 464
$ws.Range("K100").Value = 928
$ws.Range("M100").Value = -387
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954
$ws.Range("H121").Value = 67979.5
$ws.Range("J121").Value = 67979.5
$ws.Range("L121").Value = 67979.5
$ws.Range("N121").Value = -71473.5
$ws.Range("H126").Value = 2535.8823
$ws.Range("I126").Value = 1807.3334
$ws.Range("K126").Value = 5422.0002
$ws.Range("M126").Value = -2952.0002
$ws.Range("H132").Value = 41694332
$ws.Range("I132").Value = 62506744
$ws.Range("K132").Value = 187520232
$ws.Range("M132").Value = -187517702
$ws.Range("H136").Value = 83421500
$ws.Range("I136").Value = 333334660
$ws.Range("J136").Value = 117111.664
$ws.Range("K136").Value = 1000003980
$ws.Range("L136").Value = 351334.992
$ws.Range("M136").Value = -1000001430
$ws.Range("N136").Value = -356434.992
